$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and Report Covering date range) ---
$ws.Range("A8").Value = "Volume 31   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"

# --- Type-changing cells: number -> text or text -> number ---
# C14: was text "0" (shared index 20) -> becomes number 3 (style matches F/I columns, #,##0)
$c = $ws.Range("C14")
$c.NumberFormat = "#,##0"
$c.Value = 3

# C22: was text "0" -> becomes number 1
$c = $ws.Range("C22")
$c.NumberFormat = "#,##0"
$c.Value = 1

# C28: was text "0" -> becomes number 2
$c = $ws.Range("C28")
$c.NumberFormat = "#,##0"
$c.Value = 2

# G14: was number 1 -> becomes text "0" (shared string, style matches C/D/E text cells)
$dst = $ws.Range("G14")
$dst.NumberFormat = "@"
$dst.Value = "0"
$ws.Range("D14").Copy() | Out-Null
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# H14: was number 0 -> becomes text "***.*" (shared string, style matches C/D/E text cells)
$dst = $ws.Range("H14")
$dst.NumberFormat = "@"
$dst.Value = "***.*"
$ws.Range("D14").Copy() | Out-Null
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Plain numeric value updates (same cell style throughout) ---
$ws.Range("F14").Value = 4
$ws.Range("I14").Value = 6
$ws.Range("K14").Value = 500
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = -50
$ws.Range("D15").Value = 1
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -33.333333333333
$ws.Range("M15").Value = 100
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = 22.727272727272
$ws.Range("I16").Value = 64
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = -3.030303030303
$ws.Range("L16").Value = -14.666666666666
$ws.Range("M16").Value = -17.948717948717
$ws.Range("N16").Value = -87.044534412955
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 45
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 144
$ws.Range("J17").Value = 138
$ws.Range("K17").Value = 4.347826086956
$ws.Range("L17").Value = 5.109489051094
$ws.Range("M17").Value = 97.260273972602
$ws.Range("N17").Value = -22.994652406417
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 51
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = 34.210526315789
$ws.Range("L18").Value = 15.90909090909
$ws.Range("M18").Value = 18.60465116279
$ws.Range("N18").Value = -88.302752293578
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 32.5
$ws.Range("I19").Value = 143
$ws.Range("J19").Value = 110
$ws.Range("K19").Value = 30
$ws.Range("L19").Value = 11.71875
$ws.Range("M19").Value = 191.836734693878
$ws.Range("N19").Value = -20.994475138121
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -27.272727272727
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = -32.258064516129
$ws.Range("L20").Value = -32.258064516129
$ws.Range("M20").Value = 75
$ws.Range("N20").Value = -81.659388646288
$ws.Range("C21").Value = 51
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = 24.390243902439
$ws.Range("F21").Value = 167
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = 15.972222222222
$ws.Range("I21").Value = 454
$ws.Range("J21").Value = 421
$ws.Range("K21").Value = 7.838479809976
$ws.Range("L21").Value = 0.442477876106
$ws.Range("M21").Value = 67.527675276752
$ws.Range("N21").Value = -70.652876535229
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 40
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 40
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -66.666666666666
$ws.Range("L23").Value = -37.5
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 61
$ws.Range("H24").Value = 39.344262295082
$ws.Range("I24").Value = 206
$ws.Range("J24").Value = 181
$ws.Range("K24").Value = 13.812154696132
$ws.Range("L24").Value = -22.556390977443
$ws.Range("M24").Value = 1.980198019801
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -80
$ws.Range("F25").Value = 17
$ws.Range("H25").Value = 88.888888888888
$ws.Range("I25").Value = 47
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = -6
$ws.Range("L25").Value = -68.666666666666
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 46.153846153846
$ws.Range("F26").Value = 75
$ws.Range("G26").Value = 67
$ws.Range("H26").Value = 11.940298507462
$ws.Range("I26").Value = 240
$ws.Range("J26").Value = 232
$ws.Range("K26").Value = 3.448275862068
$ws.Range("L26").Value = 25
$ws.Range("M26").Value = 8.597285067873
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = -30.76923076923
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -62.5
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 18
$ws.Range("K28").Value = -11.111111111111
$ws.Range("L28").Value = -15.78947368421
$ws.Range("C29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 300
$ws.Range("I29").Value = 6
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = -25
$ws.Range("N29").Value = -75
$ws.Range("C30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 300
$ws.Range("I30").Value = 6
$ws.Range("K30").Value = 20
$ws.Range("L30").Value = -25
$ws.Range("M30").Value = -25
$ws.Range("N30").Value = -72.727272727272
